# Fruta / hortaliza, semanal
# Inserts a new weekly record for "Piña" (Vega Modelo de Temuco) at row 427,
# pushing the existing rows 427-521 down to 428-522.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 427 - shifts rows 427..521 down to 428..522
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new weekly data point.
$ws.Range("A427").Value = 10
$ws.Range("B427").Value = "Vega Modelo de Temuco"
$ws.Range("C427").Value = "La Araucanía"
$ws.Range("D427").Value = 44785
$ws.Range("E427").Value = 9
$ws.Range("F427").Value = "Fruta"
$ws.Range("G427").Value = 100108
$ws.Range("H427").Value = "Tropicales y subtropicales"
$ws.Range("I427").Value = 100108005
$ws.Range("J427").Value = "Piña"
$ws.Range("K427").Value = "Caramelo"
$ws.Range("L427").Value = "Segunda"
$ws.Range("M427").Value = 130
$ws.Range("N427").Value = 20000
$ws.Range("O427").Value = 22000
$ws.Range("P427").Value = 21154
$ws.Range("Q427").Value = "$/caja 14 unidades"
$ws.Range("R427").Value = "Ecuador"
$ws.Range("S427").Value = 1511
$ws.Range("T427").Value = 14
